# Update gh-pages output values (F column = "想去人数") to the freshly
# crawled numbers, on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value  = 5455
$wsExhibit.Range("F9").Value  = 7366
$wsExhibit.Range("F12").Value = 3761
$wsExhibit.Range("F24").Value = 5135
$wsExhibit.Range("F29").Value = 7594
$wsExhibit.Range("F32").Value = 2133
$wsExhibit.Range("F35").Value = 1159
$wsExhibit.Range("F46").Value = 1991
$wsExhibit.Range("F47").Value = 110
$wsExhibit.Range("F49").Value = 1204

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value  = 5455
$wsAll.Range("F10").Value = 3761
$wsAll.Range("F24").Value = 5135
$wsAll.Range("F29").Value = 7594
$wsAll.Range("F32").Value = 2133
$wsAll.Range("F35").Value = 1159
$wsAll.Range("F46").Value = 1991
$wsAll.Range("F47").Value = 110

$wb.Save()
